$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update view-count column F for several rows
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 1325
$wsExhibit.Range("F8").Value = 10463
$wsExhibit.Range("F11").Value = 278
$wsExhibit.Range("F13").Value = 673
$wsExhibit.Range("F14").Value = 11955
$wsExhibit.Range("F15").Value = 12361

# Sheet "全部类型" (sheet4): same events appear one row lower, update accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 1325
$wsAll.Range("F9").Value = 10463
$wsAll.Range("F12").Value = 278
$wsAll.Range("F14").Value = 673
$wsAll.Range("F15").Value = 11955
$wsAll.Range("F16").Value = 12361
